# Regenerate the "K" column (column G) values for save_data.
# The workbook stores pre-computed statistics written by an external data
# pipeline (no live formulas in the sheet), so we reproduce the effect of
# that regen by writing the new literal values directly into column G,
# matching the rows that actually changed between the two committed
# snapshots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 2
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    14 = 1
    16 = 2
    18 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
